$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "obs" unit label used across K2:K54 (shared string)
$ws.Range("K2:K54").Value = "dN/N_inel-dy-dpT"

# 2. Widen column K to fit the new longer text
$ws.Range("K1").ColumnWidth = 15.45

# 3. Update the active selection to K2:K54
$ws.Range("K2:K54").Select() | Out-Null
